$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.396.55"
$ws.Range("E2").Value = "  -2.59%  "
$ws.Range("D3").Value = "3.782.03"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.04%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "593.17"
$r.ClearFormats()
$ws.Range("E5").Value = "  -1.10%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "166.10"
$r.ClearFormats()
$ws.Range("E6").Value = "  -2.46%  "
$ws.Range("D7").Value = "3.780.62"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -1.57%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.159"
$r.ClearFormats()
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("E11").Value = "  -2.51%  "
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("E13").Value = "  -3.55%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "35.88"
$r.ClearFormats()
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Value = "4.417.27"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "3.769.66"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "67.381.35"
$ws.Range("E17").Value = "  -2.57%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "18.05"
$r.ClearFormats()
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E20").Value = "  -1.82%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "10.20"
$r.ClearFormats()
$ws.Range("E21").Value = "  -7.20%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "457.97"
$r.ClearFormats()
$ws.Range("E22").Value = "  -2.93%  "
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("E24").Value = "  +1.43%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "83.36"
$r.ClearFormats()
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("E26").Value = "  -4.96%  "
$ws.Range("E27").Value = "  -3.20%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "10.00"
$r.ClearFormats()
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("E30").Value = "  -1.74%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "29.79"
$r.ClearFormats()
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("E32").Value = "  -2.36%  "
$ws.Range("E33").Value = "  -3.96%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "9.15"
$r.ClearFormats()
$ws.Range("E34").Value = "  -2.98%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.ClearFormats()
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "3.735.99"
$ws.Range("E36").Value = "  -0.51%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "0.0995"
$r.ClearFormats()
$ws.Range("E37").Value = "  -2.75%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "3.30"
$r.ClearFormats()
$ws.Range("E38").Value = "  -6.97%  "
$ws.Range("E39").Value = "  -1.58%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.993"
$r.ClearFormats()
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("E41").Value = "  -2.99%  "
$ws.Range("E42").Value = "  -0.07%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "44.17"
$r.ClearFormats()
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("E46").Value = "  -4.52%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "8.34"
$r.ClearFormats()
$ws.Range("E47").Value = "  -3.78%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "147.96"
$r.ClearFormats()
$ws.Range("E48").Value = "  +1.14%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "392.18"
$r.ClearFormats()
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("E50").Value = "  -8.37%  "
$ws.Range("D51").Value = "2.753.32"
$ws.Range("E51").Value = "  +1.80%  "
